$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B148").Value = 65258
$ws.Range("B149").Value = 64196
$ws.Range("B195").Value = 64350
$ws.Range("E195").Value = 70.63
$ws.Range("F195").Value = 2
$ws.Range("G195").Value = 132.88
$ws.Range("B196").Value = 57756
$ws.Range("E196").Value = 79.37
$ws.Range("F196").Value = -100
$ws.Range("G196").Value = -6644
$ws.Range("F274").Value = 97
$ws.Range("G274").Value = 8641.73
$ws.Range("F276").Value = 122
$ws.Range("G276").Value = 12866.12
$ws.Range("F277").Value = 194
$ws.Range("G277").Value = 21784.26
$ws.Range("F278").Value = 102
$ws.Range("G278").Value = 7578.6
$ws.Range("F279").Value = 68
$ws.Range("G279").Value = 5052.4
$ws.Range("B280").Value = 95495.12
$ws.Range("B322").Value = 66188
$ws.Range("C322").Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("D322").Value = 315.8
$ws.Range("E322").Value = 377.31
$ws.Range("F322").Value = 29
$ws.Range("G322").Value = 9158.200000000001
$ws.Range("B323").Value = 48719
$ws.Range("C323").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D323").Value = 295.75
$ws.Range("E323").Value = 353.35
$ws.Range("F323").Value = -82
$ws.Range("G323").Value = -24251.5
$ws.Range("B366").Value = 64983
$ws.Range("C366").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F366").Value = 6
$ws.Range("G366").Value = 514.08
$ws.Range("B367").Value = 66194
$ws.Range("C367").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F367").Value = 35
$ws.Range("G367").Value = 2998.8
$ws.Range("B374").Value = 61610
$ws.Range("E374").Value = 122.71
$ws.Range("F374").Value = -58
$ws.Range("G374").Value = -5957.18
$ws.Range("B375").Value = 63565
$ws.Range("E375").Value = 109.19
$ws.Range("F375").Value = 60
$ws.Range("G375").Value = 6162.6
$ws.Range("F431").Value = 13
$ws.Range("G431").Value = 2198.82
$ws.Range("B450").Value = 88527.83
$ws.Range("F455").Value = 310
$ws.Range("G455").Value = 43582.9
$ws.Range("B457").Value = 98062.38
$ws.Range("B555").Value = 64922
$ws.Range("E555").Value = 20.98
$ws.Range("F555").Value = 0
$ws.Range("G555").Value = 0
$ws.Range("B556").Value = 45706
$ws.Range("E556").Value = 23.58
$ws.Range("F556").Value = -207
$ws.Range("G556").Value = -4084.11
$ws.Range("B563").Value = 45709
$ws.Range("E563").Value = 15.69
$ws.Range("F563").Value = -302
$ws.Range("G563").Value = -3971.3
$ws.Range("B564").Value = 64925
$ws.Range("E564").Value = 13.97
$ws.Range("F564").Value = 0
$ws.Range("G564").Value = 0
$ws.Range("B639").Value = 64810
$ws.Range("E639").Value = 291.22
$ws.Range("F639").Value = 1
$ws.Range("G639").Value = 273.92
$ws.Range("B640").Value = 53319
$ws.Range("E640").Value = 310.64
$ws.Range("F640").Value = -6
$ws.Range("G640").Value = -1643.52
$ws.Range("B658").Value = 60025
$ws.Range("E658").Value = 37.22
$ws.Range("F658").Value = -98
$ws.Range("G658").Value = -3217.34
$ws.Range("B659").Value = 64833
$ws.Range("E659").Value = 34.9
$ws.Range("F659").Value = 88
$ws.Range("G659").Value = 2889.04
$ws.Range("B668").Value = 60022
$ws.Range("E668").Value = 37.22
$ws.Range("F668").Value = -113
$ws.Range("G668").Value = -3709.79
$ws.Range("B669").Value = 64830
$ws.Range("E669").Value = 34.9
$ws.Range("F669").Value = 88
$ws.Range("G669").Value = 2889.04
$ws.Range("F803").Value = 0
$ws.Range("G803").Value = 0
$ws.Range("F804").Value = 0
$ws.Range("G804").Value = 0
$ws.Range("F809").Value = 10
$ws.Range("G809").Value = 852
$ws.Range("F812").Value = 0
$ws.Range("G812").Value = 0
$ws.Range("F816").Value = 35
$ws.Range("G816").Value = 4996.95
$ws.Range("F822").Value = 4
$ws.Range("G822").Value = 132.24
$ws.Range("F823").Value = 0
$ws.Range("G823").Value = 0
$ws.Range("F825").Value = 31
$ws.Range("G825").Value = 11712.11
$ws.Range("F826").Value = 527
$ws.Range("G826").Value = 54233.57
$ws.Range("B828").Value = 65362
$ws.Range("F828").Value = 0
$ws.Range("G828").Value = 0
$ws.Range("B829").Value = 65079
$ws.Range("F829").Value = 6
$ws.Range("G829").Value = 245.22
$ws.Range("F831").Value = 42
$ws.Range("G831").Value = 1982.82
$ws.Range("F836").Value = 0
$ws.Range("G836").Value = 0
$ws.Range("B837").Value = 258144.44
$ws.Range("B937").Value = 3606574.78
$ws.Range("B938").Value = 3606574.78
